$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data corrections (row 16) ---
$ws.Range("C16").Value = 64254

# --- Data corrections (row 19) ---
$ws.Range("D19").Value = 12324
$ws.Range("K19").Value = 514

# --- Data corrections (row 21) ---
$ws.Range("E21").Value = 786

# --- Data corrections (row 25) ---
$ws.Range("D25").Value = 94310
$ws.Range("F25").Value = 98340
$ws.Range("H25").Value = 84530
$ws.Range("J25").Value = 163370
$ws.Range("M25").Value = 34400

# --- Data corrections (row 26) ---
$ws.Range("C26").Value = 1942
$ws.Range("H26").Value = 1890945

# --- Data corrections (row 27) ---
$ws.Range("F27").Value = 92548
$ws.Range("J27").Value = 73742

# --- Data corrections (row 29) ---
$ws.Range("M29").Value = 8150

# --- Data corrections (row 30) ---
$ws.Range("F30").Value = 12280

# --- Data corrections (row 34) ---
$ws.Range("E34").Value = 27441

# --- Data corrections (row 35) ---
$ws.Range("M35").Value = 1064120

# --- Data corrections (row 36) ---
$ws.Range("C36").Value = 42750
$ws.Range("E36").Value = 84753
$ws.Range("K36").Value = 103970

# --- Data corrections (row 37) ---
$ws.Range("C37").Value = 456900

# --- Data corrections (row 40) ---
$ws.Range("N40").Value = 1020

# --- Data corrections (row 41) ---
$ws.Range("I41").Value = 38920

# --- Row height fixes: rows 33 and 53 had an explicit custom row height (16pt);
# restore them to the sheet's standard height like the surrounding rows.
$ws.Rows.Item(33).AutoFit()
$ws.Rows.Item(53).AutoFit()

# --- View state: zoom + scroll position + selection, as last left by the editor ---
$win = $excel.ActiveWindow
$win.Zoom = 116
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("F42").Select()
